$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to stay text for numeric-looking values,
# then restore default styling so no style/format diff is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.032.61"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.965.49"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "591.61"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "142.15"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.959.29"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").Value = "5.99"
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "33.86"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "3.459.02"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "61.052.16"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "6.86"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").Value = "2.969.86"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "445.75"
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("D21").Value = "13.91"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").Value = "81.39"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "10.46"
$ws.Range("E25").Value = "  +6.05%  "
$ws.Range("D26").Value = "2.16"
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("D27").Value = "11.90"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "2.66"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "2.03"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").Value = "26.95"
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "0.0₃0804"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "5.73"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "49.96"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "8.91"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.120"
$ws.Range("E41").Value = "  +7.30%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.81"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").Value = "382.86"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").Value = "0.0347"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "0.266"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").Value = "38.03"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "2.676.54"
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").Value = "130.01"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D50").Value = "0.106"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").Value = "  -1.24%  "

$ws.Range("D2:D51").Style = "Normal"

Write-Output "done"
